$d = $word.ActiveDocument

# 1. Replace the paragraph text with the documentation placeholder.
$old = "Alle validierungen sind mit plain javascript geschreiben und findet man im beiliegendem main.js file. Das File findet man unter webroot -> js."
$new = "< insert documentation >"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# 2. Remove the now-stray trailing empty paragraph that followed it.
#    (It is the last paragraph in the document, so its own paragraph mark
#    can't be deleted on its own - extend the range to swallow the
#    preceding paragraph mark instead, merging the two paragraphs.)
$paras = $d.Paragraphs
$count = $paras.Count
$last = $paras.Item($count)
if ($last.Range.Text.Trim() -eq "") {
    $prev = $paras.Item($count - 1)
    $r = $d.Range($prev.Range.End - 1, $last.Range.End)
    $r.Delete()
}
